$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure date text (A13)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-10
$ws.Range("D2").Value = 0.1001802189684367
$ws.Range("E2").Value = 0.02620239897519516

$ws.Range("D3").Value = 0.1095316566013352
$ws.Range("E3").Value = 0.008483501053124343

$ws.Range("D4").Value = 0.1182090664158832
$ws.Range("E4").Value = 0.009428729894620291

$ws.Range("D5").Value = 0.1378192271621025
$ws.Range("E5").Value = 0.01245133536882959

$ws.Range("D6").Value = 0.1364134042850748
$ws.Range("E6").Value = 0.004294223576672529

$ws.Range("D7").Value = 0.1419047099070783
$ws.Range("E7").Value = 0.01672989008851267

$ws.Range("D8").Value = 0.127135251978111
$ws.Range("E8").Value = 0.01753864447086806

$ws.Range("D9").Value = 0.1288064646819783
$ws.Range("E9").Value = 0.009987566250788671

$ws.Range("E10").Value = 0.01286085170730233

# Restore sheet protection (the sheet was protected before this edit;
# the original password is unknown so we re-apply protection without one)
$ws.Protect()
